$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2446
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 2500
$ws.Range("M51").Value = -2016
$ws.Range("H129").Value = 1134.0186
$ws.Range("I129").Value = 460
$ws.Range("J129").Value = 1218.2709
$ws.Range("K129").Value = 1380
$ws.Range("L129").Value = 3654.8127
$ws.Range("M129").Value = 3620
$ws.Range("N129").Value = -13654.8127
$ws.Range("H132").Value = 1979.36
$ws.Range("I132").Value = 2098.848
$ws.Range("J132").Value = 605.25
$ws.Range("K132").Value = 6296.544
$ws.Range("L132").Value = 1815.75
$ws.Range("M132").Value = -3766.544
$ws.Range("N132").Value = -6875.75
$ws.Range("H141").Value = 1651.95
$ws.Range("I141").Value = 1391.1111
$ws.Range("K141").Value = 4173.3333
$ws.Range("M141").Value = 1006.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1694.8438
$ws.Range("I2").Value = 1571.2222
$ws.Range("K2").Value = 1571.2222
$ws.Range("M2").Value = -1458.2222
$ws.Range("H45").Value = 3942
$ws.Range("I45").Value = 4440.375
$ws.Range("J45").Value = 3543.3
$ws.Range("K45").Value = 4440.375
$ws.Range("L45").Value = 3543.3
$ws.Range("M45").Value = -4063.375
$ws.Range("N45").Value = -4297.3
$ws.Range("H63").Value = 3156.4
$ws.Range("I63").Value = 2600
$ws.Range("J63").Value = 3712.8
$ws.Range("K63").Value = 2600
$ws.Range("L63").Value = 3712.8
$ws.Range("M63").Value = -1914
$ws.Range("N63").Value = -5084.8
$ws.Range("H66").Value = 3156.4
$ws.Range("I66").Value = 2600
$ws.Range("J66").Value = 3712.8
$ws.Range("K66").Value = 13000
$ws.Range("L66").Value = 18564
$ws.Range("M66").Value = -9568
$ws.Range("N66").Value = -25428
$ws.Range("H116").Value = 1694.8438
$ws.Range("I116").Value = 1571.2222
$ws.Range("K116").Value = 1571.2222
$ws.Range("M116").Value = 722.7778000000001
$ws.Range("H132").Value = 20582.555
$ws.Range("I132").Value = 1901.3478
$ws.Range("J132").Value = 127999.5
$ws.Range("K132").Value = 5704.0434
$ws.Range("L132").Value = 383998.5
$ws.Range("M132").Value = -3174.0434
$ws.Range("N132").Value = -389058.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1694.8438
$ws.Range("I3").Value = 1571.2222
$ws.Range("K3").Value = 1571.2222
$ws.Range("M3").Value = -1457.2222
$ws.Range("H86").Value = 1664.4062
$ws.Range("I86").Value = 1432.8846
$ws.Range("J86").Value = 2667.6667
$ws.Range("K86").Value = 1432.8846
$ws.Range("L86").Value = 2667.6667
$ws.Range("M86").Value = -309.8846000000001
$ws.Range("N86").Value = -4913.6667
$ws.Range("H89").Value = 1664.4062
$ws.Range("I89").Value = 1432.8846
$ws.Range("J89").Value = 2667.6667
$ws.Range("K89").Value = 7164.423000000001
$ws.Range("L89").Value = 13338.3335
$ws.Range("M89").Value = -1548.423000000001
$ws.Range("N89").Value = -24570.3335
$ws.Range("H105").Value = 3401.25
$ws.Range("I105").Value = 3310.9092
$ws.Range("J105").Value = 3600
$ws.Range("K105").Value = 3310.9092
$ws.Range("L105").Value = 3600
$ws.Range("M105").Value = -1563.9092
$ws.Range("N105").Value = -7094
$ws.Range("H134").Value = 2679.7114
$ws.Range("I134").Value = 2844.4773
$ws.Range("K134").Value = 8533.4319
$ws.Range("M134").Value = -5998.4319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 47621790
$ws.Range("I62").Value = 50002430
$ws.Range("J62").Value = 9006
$ws.Range("K62").Value = 50002430
$ws.Range("L62").Value = 9006
$ws.Range("M62").Value = -50001806
$ws.Range("N62").Value = -10254
$ws.Range("H65").Value = 47621790
$ws.Range("I65").Value = 50002430
$ws.Range("J65").Value = 9006
$ws.Range("K65").Value = 250012150
$ws.Range("L65").Value = 45030
$ws.Range("M65").Value = -250009030
$ws.Range("N65").Value = -51270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 927.55554
$ws.Range("I23").Value = 537.5
$ws.Range("J23").Value = 1039
$ws.Range("K23").Value = 1612.5
$ws.Range("L23").Value = 3117
$ws.Range("M23").Value = -1377.5
$ws.Range("N23").Value = -3587
$ws.Range("H64").Value = 6900
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 6900
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 20700
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -21240
$ws.Range("H67").Value = 6900
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 6900
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 20700
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -22572
$ws.Range("H76").Value = 4754.5
$ws.Range("J76").Value = 4754.5
$ws.Range("L76").Value = 14263.5
$ws.Range("N76").Value = -15029.5
$ws.Range("H79").Value = 4754.5
$ws.Range("J79").Value = 4754.5
$ws.Range("L79").Value = 14263.5
$ws.Range("N79").Value = -16915.5
$ws.Range("H94").Value = 4025.3809
$ws.Range("I94").Value = 700
$ws.Range("J94").Value = 4191.65
$ws.Range("K94").Value = 2100
$ws.Range("L94").Value = 12574.95
$ws.Range("M94").Value = -1424
$ws.Range("N94").Value = -13926.95
$ws.Range("H100").Value = 6200
$ws.Range("J100").Value = 6200
$ws.Range("L100").Value = 18600
$ws.Range("N100").Value = -20222
$ws.Range("H109").Value = 6239.273
$ws.Range("I109").Value = 737
$ws.Range("J109").Value = 6789.5
$ws.Range("K109").Value = 2211
$ws.Range("L109").Value = 20368.5
$ws.Range("M109").Value = -1171
$ws.Range("N109").Value = -22448.5
$ws.Range("H112").Value = 1957.5
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 6030
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 18090
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -20306
$ws.Range("H115").Value = 5713.5557
$ws.Range("I115").Value = 3028
$ws.Range("J115").Value = 6049.25
$ws.Range("K115").Value = 9084
$ws.Range("L115").Value = 18147.75
$ws.Range("M115").Value = -7909
$ws.Range("N115").Value = -20497.75
$ws.Range("H120").Value = 15606
$ws.Range("I120").Value = 9015
$ws.Range("K120").Value = 27045
$ws.Range("M120").Value = -22207
$ws.Range("H121").Value = 9265405
$ws.Range("I121").Value = 505
$ws.Range("J121").Value = 11118385
$ws.Range("K121").Value = 1515
$ws.Range("L121").Value = 33355155
$ws.Range("M121").Value = -205
$ws.Range("N121").Value = -33357775
$ws.Range("H122").Value = 457.04166
$ws.Range("I122").Value = 237.27777
$ws.Range("K122").Value = 2135.49993
$ws.Range("M122").Value = 314.5000700000001
$ws.Range("H123").Value = 4264.6
$ws.Range("I123").Value = 2900
$ws.Range("J123").Value = 4605.75
$ws.Range("K123").Value = 8700
$ws.Range("L123").Value = 13817.25
$ws.Range("M123").Value = -6250
$ws.Range("N123").Value = -18717.25
$ws.Range("H124").Value = 1000
$ws.Range("I124").Value = 1000
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 3000
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 1910
$ws.Range("N124").ClearContents()
$ws.Range("H131").Value = 785.83
$ws.Range("J131").Value = 795.4681
$ws.Range("L131").Value = 2386.4043
$ws.Range("N131").Value = -12466.4043

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3152940.5
$ws.Range("I11").Value = 3461538.5
$ws.Range("K11").Value = 3461538.5
$ws.Range("M11").Value = -3461399.5
$ws.Range("H126").Value = 4954.357
$ws.Range("I126").Value = 3801.158
$ws.Range("J126").Value = 7388.8887
$ws.Range("K126").Value = 11403.474
$ws.Range("L126").Value = 22166.6661
$ws.Range("M126").Value = -8933.474
$ws.Range("N126").Value = -27106.6661
$ws.Range("H132").Value = 52728.273
$ws.Range("I132").Value = 10335.667
$ws.Range("K132").Value = 31007.001
$ws.Range("M132").Value = -28477.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3081.9092
$ws.Range("I132").Value = 2351
$ws.Range("K132").Value = 7053
$ws.Range("L132").Value = 10498.7145
$ws.Range("M132").Value = -4523
$ws.Range("N132").Value = -15558.7145

